$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Column G: area per segment
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G15").Formula = "=(D4-D3)*B4/100"

# Column H: total area
$ws.Range("H2").Formula = "=SUM(G2:G11)"

$ws.Range("F2").Select()
